$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "初始"
$ws.Range("B6").Value = "皮皮"
$ws.Range("B7").Value = "水管工"
$ws.Range("A8").Value = $null
$ws.Range("B8").Value = $null

$ws.Range("B6").Select()
